$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A21").Value = "How many curve shades can I create?"
$ws.Range("B21").Value = "llama3.2:latest"
$ws.Range("C21").Value = "According to the document, you can create 250 curve shades."
